$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a blank row at row 20. This pushes the old signature rows
#    (24, 25) down to (25, 26) without touching row 19's own content/style.
$ws.Rows.Item(20).Insert()

# 2. Populate the new row 20 with an exact copy (value + style) of what is
#    currently row 19 (the last "account period" data row), cell by cell so
#    the original style index for each column is preserved verbatim.
$cols = @("B","C","D","E","F","G","H","I","J")
foreach ($col in $cols) {
  $ws.Range($col + "19").Copy($ws.Range($col + "20"))
}

# 3. Row 19 is no longer the last row of the table, so give it the same
#    ("middle of table") formatting as row 18, keeping row 19's own values.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4. Refresh the "Periodo Mora" values for the 5 data rows: a new period
#    (2507) is added at the top and the oldest one (2503) now lives in the
#    newly created last row.
$ws.Range("E16").Value = "2507"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2505"
$ws.Range("E19").Value = "2504"
$ws.Range("E20").Value = "2503"

# 5. Update the account-summary figures to match the refreshed data.
$ws.Range("E11").Value = 284700
$ws.Range("F13").Value = 5
